# Auto update Excel log
#
# Appends new mmWave sensor-log rows to three worksheets (Date, Timestamp,
# Hour, Location, Value, Status columns). The previous row is copied down
# first so the new rows inherit the existing (default/general) formatting
# instead of letting Excel's autodetect turn the date/time-looking text in
# columns A-C into real date/time serial values; only the cells whose
# content actually differs from the template row are then overwritten.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$templateRow,
        [int]$row,
        [string]$timestamp,
        $value
    )

    $ws.Range("A" + $templateRow + ":F" + $templateRow).Copy()
    $ws.Range("A" + $row + ":F" + $row).PasteSpecial()

    $ws.Range("B" + $row).Value = $timestamp

    if ($value -ne $null) {
        $ws.Range("E" + $row).Value = $value
    }
}

# --- mmWave(InBed): append rows 90-94 (Value column stays "In Bed") ---
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")

Add-LogRow $wsInBed 89 90 "20:31:58" $null
Add-LogRow $wsInBed 89 91 "20:31:59" $null
Add-LogRow $wsInBed 89 92 "20:32:01" $null
Add-LogRow $wsInBed 89 93 "20:32:03" $null
Add-LogRow $wsInBed 89 94 "20:32:04" $null

# --- mmWave(BR): append rows 86-88 (numeric Value column) ---
$wsBR = $wb.Worksheets.Item("mmWave(BR)")

Add-LogRow $wsBR 85 86 "20:32:00" 7
Add-LogRow $wsBR 85 87 "20:32:02" 6
Add-LogRow $wsBR 85 88 "20:32:04" 2

# --- mmWave(HR): append rows 86-88 (numeric Value column) ---
$wsHR = $wb.Worksheets.Item("mmWave(HR)")

Add-LogRow $wsHR 85 86 "20:32:00" 55
Add-LogRow $wsHR 85 87 "20:32:01" 54
Add-LogRow $wsHR 85 88 "20:32:03" 50
